# Update "想去人数" (interested-count) figures in the 北京-漫展信息 workbook
# Target values per commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 28   # was 26
$ws.Range("F4").Value = 5980   # was 5976
$ws.Range("F5").Value = 71   # was 70
$ws.Range("F6").Value = 3032   # was 3028
$ws.Range("F7").Value = 1297   # was 1296
$ws.Range("F9").Value = 108   # was 107
$ws.Range("F11").Value = 38   # was 36
$ws.Range("F12").Value = 333   # was 328
$ws.Range("F13").Value = 4468   # was 4459
$ws.Range("F14").Value = 4468   # was 4459
$ws.Range("F21").Value = 6932   # was 6916
$ws.Range("F22").Value = 239   # was 237
$ws.Range("F23").Value = 115   # was 112
$ws.Range("F24").Value = 480   # was 479
$ws.Range("F25").Value = 1277   # was 1275
$ws.Range("F27").Value = 1651   # was 1650
$ws.Range("F28").Value = 16   # was 15
$ws.Range("F30").Value = 6048   # was 6045
$ws.Range("F35").Value = 438   # was 436
$ws.Range("F36").Value = 6128   # was 6116
$ws.Range("F42").Value = 2421   # was 2420
$ws.Range("F43").Value = 26   # was 25
$ws.Range("F45").Value = 1011   # was 1010
$ws.Range("F47").Value = 366   # was 361
$ws.Range("F49").Value = 21   # was 20

# --- Sheet "本地生活" (local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1423   # was 1421

# --- Sheet "全部类型" (all types, aggregate view) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1423   # was 1421
$ws.Range("F3").Value = 5980   # was 5976
$ws.Range("F4").Value = 5980   # was 5976
$ws.Range("F5").Value = 3032   # was 3028
$ws.Range("F6").Value = 1297   # was 1296
$ws.Range("F10").Value = 38   # was 36
$ws.Range("F12").Value = 333   # was 328
$ws.Range("F13").Value = 4468   # was 4459
$ws.Range("F14").Value = 4468   # was 4459
$ws.Range("F21").Value = 6932   # was 6916
$ws.Range("F22").Value = 239   # was 237
$ws.Range("F23").Value = 115   # was 112
$ws.Range("F24").Value = 480   # was 479
$ws.Range("F25").Value = 1277   # was 1275
$ws.Range("F28").Value = 1651   # was 1650
$ws.Range("F32").Value = 6048   # was 6045
$ws.Range("F37").Value = 438   # was 436
$ws.Range("F38").Value = 6128   # was 6116
$ws.Range("F44").Value = 2421   # was 2420
$ws.Range("F45").Value = 26   # was 25
$ws.Range("F46").Value = 1011   # was 1010
$ws.Range("F48").Value = 366   # was 361
$ws.Range("F50").Value = 21   # was 20

$wb.Save()
